$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 1810
$ws.Range("E2").Value = 46200602503
$ws.Range("X2").Value = "DN4127460129014"

# Row 3
$ws.Range("A3").Value = 1811
$ws.Range("E3").Value = 46200602504
$ws.Range("X3").Value = "DN4127460129015"
